$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.953.77'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '2.214.84'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  -0.05%  '
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '289.19'
$c.Style = $s
$ws.Range("E5").Value = '  -1.36%  '
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '87.36'
$c.Style = $s
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +0.34%  '
$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '30.36'
$c.Style = $s
$ws.Range("E10").Value = '  +0.91%  '
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0775'
$c.Style = $s
$ws.Range("E11").Value = '  -2.59%  '
$ws.Range("E12").Value = '  +2.62%  '
$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.45'
$c.Style = $s
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D14").Value = '2.553.86'
$ws.Range("E14").Value = '  -0.86%  '
$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '13.93'
$c.Style = $s
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = '2.199.67'
$ws.Range("E16").Value = '  -1.46%  '
$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.727'
$c.Style = $s
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '39.892.06'
$ws.Range("E18").Value = '  +0.23%  '
$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.47'
$c.Style = $s
$ws.Range("E19").Value = '  +9.20%  '
$ws.Range("D20").Value = '0.0₃0881'
$ws.Range("E20").Value = '  -0.99%  '
$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.79'
$c.Style = $s
$ws.Range("E21").Value = '  +0.20%  '
$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '65.53'
$c.Style = $s
$ws.Range("E22").Value = '  +0.27%  '
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '236.84'
$c.Style = $s
$ws.Range("E23").Value = '  +2.01%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +1.13%  '
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.82'
$c.Style = $s
$ws.Range("E26").Value = '  -1.13%  '
$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '22.48'
$c.Style = $s
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("E28").Value = '  -0.21%  '
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.19'
$c.Style = $s
$ws.Range("E29").Value = '  -0.23%  '
$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '155.84'
$c.Style = $s
$ws.Range("E30").Value = '  +0.84%  '
$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '31.64'
$c.Style = $s
$ws.Range("E31").Value = '  -3.43%  '
$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = $s
$ws.Range("E32").Value = '  -0.10%  '
$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.92'
$c.Style = $s
$ws.Range("E33").Value = '  +1.68%  '
$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0715'
$c.Style = $s
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("E36").Value = '  +6.52%  '
$ws.Range("E37").Value = '  -0.19%  '
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '15.70'
$c.Style = $s
$ws.Range("E38").Value = '  -4.11%  '
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("D41").Value = '2.104.60'
$ws.Range("E41").Value = '  +7.94%  '
$ws.Range("E42").Value = '  +2.45%  '
$ws.Range("E43").Value = '  -1.32%  '
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.95'
$c.Style = $s
$ws.Range("E44").Value = '  +6.26%  '
$ws.Range("E45").Value = '  -1.14%  '
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '17.40'
$c.Style = $s
$ws.Range("E46").Value = '  +6.47%  '
$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.65'
$c.Style = $s
$ws.Range("E47").Value = '  +1.83%  '
$ws.Range("D48").Value = '2.428.33'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '68.99'
$c.Style = $s
$ws.Range("E49").Value = '  -2.46%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.44'
$c.Style = $s
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '88.40'
$c.Style = $s
$ws.Range("E51").Value = '  -0.39%  '